$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------------
# 1. Insert the new rows needed.
#    Before: row75/76 = colMeans/apply (matrix rows), row77/78 = blank.
#    After : row75 = new "sd()" vector row, row76/77 = colMeans/apply (shifted),
#            row78/79 = new "qnorm/qt" distribution rows, row80/81 = blank (shifted),
#            row82..87 = new blank rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(75).Insert()

$ws.Rows.Item(78).Insert()
$ws.Rows.Item(78).Insert()

$ws.Rows.Item(82).Insert()
$ws.Rows.Item(82).Insert()
$ws.Rows.Item(82).Insert()
$ws.Rows.Item(82).Insert()
$ws.Rows.Item(82).Insert()
$ws.Rows.Item(82).Insert()

# ---------------------------------------------------------------------------
# 2. Fix up formatting on newly inserted rows by copying formats from
#    neighbouring rows that already carry the correct style.
# ---------------------------------------------------------------------------
$ws.Range("A74:E74").Copy()
$ws.Range("A75:E75").PasteSpecial(-4122)

$ws.Range("A76:E76").Copy()
$ws.Range("A78:E79").PasteSpecial(-4122)

$ws.Range("A80:E80").Copy()
$ws.Range("A82:E87").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Row heights for the new content rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(75).RowHeight = 48
$ws.Rows.Item(78).RowHeight = 32
$ws.Rows.Item(79).RowHeight = 64

# ---------------------------------------------------------------------------
# 4. Row 75 content: 向量 / sd() / (reuse var() explanation) / 标准差
# ---------------------------------------------------------------------------
$ws.Range("A75").Value = "向量"
$ws.Range("C75").Value = "sd(x[, na.rm = TRUE/FALSE])"
$ws.Range("D74").Copy()
$ws.Range("D75").PasteSpecial(-4104)
$excel.CutCopyMode = 0
$ws.Range("E75").Value = "标准差"

# ---------------------------------------------------------------------------
# 5. Row 78 content: 分布 / qnorm(x) / explanation / description
# ---------------------------------------------------------------------------
$ws.Range("A78").Value = "分布"
$ws.Range("C78").Value = "qnorm(x)"

$d78 = $ws.Range("D78")
$d78.Value = "x: 置信水平的对应数值（99%→x=0.995，95%→x=0.975，90%→x=0.95）"
$d78.Characters(4, 9).Font.Name = "等线"
$d78.Characters(13, 1).Font.Name = "宋体"
$d78.Characters(14, 11).Font.Name = "Cascadia Code"
$d78.Characters(25, 1).Font.Name = "宋体"
$d78.Characters(26, 3).Font.Name = "Cascadia Code"
$d78.Characters(29, 1).Font.Name = "宋体"
$d78.Characters(30, 7).Font.Name = "Cascadia Code"
$d78.Characters(37, 1).Font.Name = "宋体"
$d78.Characters(38, 3).Font.Name = "Cascadia Code"
$d78.Characters(41, 1).Font.Name = "宋体"
$d78.Characters(42, 6).Font.Name = "Cascadia Code"
$d78.Characters(48, 1).Font.Name = "宋体"

$ws.Range("E78").Value = "置信水平在正态分布的对应临界值"

# ---------------------------------------------------------------------------
# 6. Row 79 content: 分布 / qt(x, df = y) / explanation / description
# ---------------------------------------------------------------------------
$ws.Range("A79").Value = "分布"
$ws.Range("C79").Value = "qt(x, df = y)"

$d79 = $ws.Range("D79")
$d79.Value = "x: 置信水平的对应数值（99%→x=0.995，95%→x=0.975，90%→x=0.95）" + $nl + "y: 自由度；自由度为n-1时，此处填写length(n-1)"
$d79.Characters(4, 10).Font.Name = "等线"
$d79.Characters(14, 11).Font.Name = "Cascadia Code"
$d79.Characters(25, 1).Font.Name = "等线"
$d79.Characters(26, 11).Font.Name = "Cascadia Code"
$d79.Characters(37, 1).Font.Name = "等线"
$d79.Characters(38, 10).Font.Name = "Cascadia Code"
$d79.Characters(48, 2).Font.Name = "等线"
$d79.Characters(50, 3).Font.Name = "Cascadia Code"
$d79.Characters(53, 8).Font.Name = "宋体"
$d79.Characters(61, 3).Font.Name = "Cascadia Code"
$d79.Characters(64, 6).Font.Name = "宋体"
$d79.Characters(70, 11).Font.Name = "Cascadia Code"

$e79 = $ws.Range("E79")
$e79.Value = "置信水平在学生t分布的对应临界值"
$e79.Characters(1, 7).Font.Name = "等线"
$e79.Characters(8, 1).Font.Name = "Cascadia Code"
$e79.Characters(9, 8).Font.Name = "等线"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 7. Sheet view / selection to match final saved state.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Range("E81").Select()
